# Rebuild Sheet1 as a tidy 3-column table (Section / Title / Page)
# out of the old wide "table of contents" layout (16 columns x 3 rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old wide layout's contents but keep per-cell formatting around so
# the header style (bold + border, currently only on column A) survives and
# can be reused/extended below; kill the leftover formatting on A2:A3 since
# those rows are being repurposed as plain data rows.
$ws.UsedRange.ClearContents()
$ws.Range("A2:A3").ClearFormats()

# ---- Header row (row 1) ------------------------------------------------
# A1 and C1 look numeric ("1.3" / "11") so a plain .Value assignment would
# silently coerce them to numbers; force text by entering them as a quoted
# formula and then collapsing the formula down to a static value in place.
$ws.Range("A1").Formula = "=""1.3"""
$ws.Range("A1").Copy()
$ws.Range("A1").PasteSpecial(-4163)

$ws.Range("B1").Value = "Ethics and integrity"

$ws.Range("C1").Formula = "=""11"""
$ws.Range("C1").Copy()
$ws.Range("C1").PasteSpecial(-4163)

# Apply the bold + bordered header style (already on A1) across the whole
# header row.
$ws.Range("A1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)

# ---- Data rows (rows 2-17) ---------------------------------------------
$data = @(
    @(1.4, "Corporate governance", 11),
    @(1.5, "Stakeholder engagement", 11),
    @(1.6, "Reporting practice", 14),
    @(2, "GRI 200: Economic standards", 17),
    @(2.1, "GRI 201: Economic performance 2016", 17),
    @(2.2, "GRI 203: Indirect economic impacts 2016", 21),
    @(2.3, "GRI 205: Anti-corruption 2016", 22),
    @(2.4, "GRI 206: Anti-competitive behavior 2016", 24),
    @(3, "GRI 300 Environmental standards", 25),
    @(3.1, "GRI 301: Materials 2016", 25),
    @(3.2, "GRI 302: Energy 2016", 26),
    @(3.3, "GRI 303: Water and effluents 2018", 27),
    @(3.4, "GRI 305: Emissions 2016", 27),
    @(3.5, "GRI 306: Effluents and waste 2016", 29),
    @(3.6, "GRI 307: Environmental compliance 2016", 30),
    @(3.7, "GRI 308: Supplier environmental assessment 2016", 30)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
